$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.585.24"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.620.88"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.64"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.089.86"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.402.21"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.636.33"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("E18").Value = "  +7.68%  "
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.02"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.73"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.32"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.73"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +10.56%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.23"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.68"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "568.75"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.26"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.28%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0844"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.76"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.25"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.48"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "167.68"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.91"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  +4.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0600"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.628"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("E49").Value = "  +5.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0963"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.39"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.55%  "
